# Update BOC USD rates (auto)
# Appends a newly captured BOC USD publish record (row 18) to the
# "All Published Values" sheet, extends the sheet's auto filter /
# filter-database defined name to cover it, and refreshes the
# dependent daily "publishes" count on the "Daily Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsAll = $wb.Worksheets.Item("All Published Values")
$wsSummary = $wb.Worksheets.Item("Daily Summary")

$targetRow = 18
$templateRow = 17

# Force the new row's cells to text formatting first so that the
# date-like / numeric-like strings are stored as plain text (matching
# the rest of the sheet) instead of being auto-converted to real
# dates or numbers.
$newRange = $wsAll.Range("A" + $targetRow + ":J" + $targetRow)
$newRange.NumberFormat = "@"

$wsAll.Cells.Item($targetRow, 1).Value = "2026-01-02"
$wsAll.Cells.Item($targetRow, 2).Value = "2026-01-02 22:55:06"
$wsAll.Cells.Item($targetRow, 3).Value = "697.85"
$wsAll.Cells.Item($targetRow, 4).Value = "697.85"
$wsAll.Cells.Item($targetRow, 5).Value = "700.79"
$wsAll.Cells.Item($targetRow, 6).Value = "700.79"
$wsAll.Cells.Item($targetRow, 7).Value = "702.88"
$wsAll.Cells.Item($targetRow, 8).Value = "2026/01/02 22:55:06"
$wsAll.Cells.Item($targetRow, 9).Value = "2026-01-02 14:57:20"
$wsAll.Cells.Item($targetRow, 10).Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Restore the plain/default cell style (matching the rest of the
# data rows) now that the values are safely stored as text.
for ($col = 1; $col -le 10; $col++) {
    $wsAll.Cells.Item($targetRow, $col).Style = $wsAll.Cells.Item($templateRow, $col).Style
}

# Re-apply the autofilter so its range grows to include the new row.
$wsAll.AutoFilterMode = $false
$wsAll.Range("A1:J18").AutoFilter()

# Keep the hidden _FilterDatabase defined name for this sheet in sync
# with the new autofilter range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$18"
    }
}

# Update the "publishes" count on the Daily Summary sheet (16 -> 17).
$wsSummary.Cells.Item(4, 2).Value = 17
